# Daily "cryptos list" refresh (GitHub Actions bot): updates the Price (D)
# and Volume(1h) (E) columns for most coins, and for three rows (13-15)
# the whole row's rank position changed (Polkadot/Litecoin/WrappedEther
# got reordered), so B/C/D/E all change there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E hold plain text in the source data (e.g. "244.74", "  +0.35%  "),
# not numbers. Pre-format the touched cells as Text so Excel doesn't
# auto-coerce numeric-looking strings ("244.74", "5.467", ...) into
# actual numbers when we assign .Value.
$textCells = @("D2", "E2", "D3", "E3", "E4", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "E23", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.936.20"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "1.895.30"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").Value = "244.74"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.3142"
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("D9").Value = "25.83"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("E10").Value = "  +1.16%  "
$ws.Range("D11").Value = "0.08913"
$ws.Range("E11").Value = "  +9.88%  "
$ws.Range("D12").Value = "0.7743"
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.467"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "94.86"
$ws.Range("E14").Value = "  +2.41%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.821.37"
$ws.Range("E15").Value = "  -4.83%  "
$ws.Range("D16").Value = "6.210"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "29.962.87"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "14.00"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "246.89"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "0.000007906"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("D21").Value = "2.164.24"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").Value = "8.189"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "0.1589"
$ws.Range("E25").Value = "  -4.89%  "
$ws.Range("D26").Value = "9.560"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").Value = "162.95"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").Value = "18.87"
$ws.Range("E28").Value = "  +0.59%  "
$ws.Range("D29").Value = "2.052"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").Value = "1.550"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").Value = "4.527"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "4.122"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "0.05530"
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").Value = "1.249"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("D36").Value = "0.7564"
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "2.719"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("D39").Value = "0.01969"
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("D40").Value = "2.793"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "0.4526"
$ws.Range("E41").Value = "  +2.18%  "
$ws.Range("D42").Value = "74.09"
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").Value = "6.093"
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("D44").Value = "1.093.61"
$ws.Range("E44").Value = "  -6.79%  "
$ws.Range("D45").Value = "0.8569"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "1.898"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "102.98"
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("D49").Value = "7.637"
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("D50").Value = "9.913"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("D51").Value = "2.998"
$ws.Range("E51").Value = "  -0.21%  "

# Restore the default (unstyled) cell style now that the values are
# safely stored as text, so we don't leave a stray explicit "Text"
# number format behind on cells that originally had none.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
